# Update a handful of BPC_Login step values on Sheet1:
#   - "Sendkeys"    -> "SendKeys"     (rows 2 and 3, column C / Action)
#   - "ls_usernam"  -> "ls_username"  (row 2, column G / IdentifyAttribute)
#   - "GetAttribute"-> "GetText"      (row 5, column C / Action)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "SendKeys"
$ws.Range("G2").Value = "ls_username"
$ws.Range("C3").Value = "SendKeys"
$ws.Range("C5").Value = "GetText"
